$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38; this shifts the existing rows 38:99
# down to 39:100 (and grows the used range to A1:R100).
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new weekly price record.
$ws.Range("A38").Value = 3
$ws.Range("B38").Value = "Femacal de La Calera"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44546
$ws.Range("E38").Value = 5
$ws.Range("F38").Value = 100112052
$ws.Range("G38").Value = "Albahaca"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 125
$ws.Range("K38").Value = 4000
$ws.Range("L38").Value = 4500
$ws.Range("M38").Value = 4260
$ws.Range("N38").Value = "$/docena de matas"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 710
$ws.Range("Q38").Value = 6
$ws.Range("R38").Value = "Hortaliza"

# Preserve the date number format on the new row's Fecha cell, matching
# the style already used by column D in the rest of the sheet.
$ws.Range("D38").NumberFormat = $ws.Range("D39").NumberFormat
